# Added PersonStateFingerprintIdentification to SSPs.
#
# The "Request" sheet holds a table of element mappings, each entry
# described by three consecutive rows/columns: Business Element Name,
# NEIM 3.0 Mapping, Element Description. A new entry -
# "Person State Fingerprint ID" - is inserted as row 9 (just above the
# existing "Person FBI Identification ID" row), pushing the remaining
# row down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new entry by inserting a row above the current
# row 9 ("Person FBI Identification ID" / PostConsolidation block).
# The insert carries the surrounding row formatting down automatically,
# matching the look of the other identification rows (5/6/8).
$ws.Rows("9:9").Insert() | Out-Null

# Populate the new row with the PersonStateFingerprintIdentification
# mapping data.
$ws.Range("A9").Value = "Person State Fingerprint ID"
$ws.Range("B9").Value = "An identification of a person based on a Fingerprint ID."
$ws.Range("C9").Value = "/CHcr-doc:CriminalHistoryConsolidationReport/nc:Person/CHcr-ext:PostConsolidationIdentifiers/j:PersonStateFingerprintIdentification/nc:IdentificationID"

# Leave the cursor where the author left it after adding the row.
$ws.Range("B15").Select() | Out-Null
